$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 4148
$ws.Range("K3").Value = 4247
$ws.Range("K4").Value = 856
$ws.Range("K5").Value = 307
$ws.Range("K6").Value = 4750
$ws.Range("K7").Value = 14308

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("K2").Value = 47
$ws.Range("K7").Value = 193

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K2").Value = 273
$ws.Range("K3").Value = 289
$ws.Range("K6").Value = 321
$ws.Range("K7").Value = 964

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("K2").Value = 106
$ws.Range("K3").Value = 105
$ws.Range("K7").Value = 304

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K2").Value = 162
$ws.Range("K3").Value = 224
$ws.Range("K5").Value = 12
$ws.Range("K6").Value = 171
$ws.Range("K7").Value = 596

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K3").Value = 161
$ws.Range("K7").Value = 484

$ws = $wb.Worksheets.Item('New City')
$ws.Range("K2").Value = 97
$ws.Range("K3").Value = 84
$ws.Range("K7").Value = 330

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K7").Value = 420
$ws.Range("K8").Value = 964
$ws.Range("K11").Value = 280
$ws.Range("K14").Value = 80
$ws.Range("K15").Value = 147
$ws.Range("K16").Value = 43
$ws.Range("K20").Value = 321
$ws.Range("K23").Value = 146
$ws.Range("K25").Value = 66
$ws.Range("K29").Value = 751
$ws.Range("K33").Value = 596
$ws.Range("K36").Value = 181
$ws.Range("K37").Value = 484
$ws.Range("K42").Value = 521
$ws.Range("K44").Value = 131
$ws.Range("K45").Value = 17
$ws.Range("K46").Value = 33
$ws.Range("K48").Value = 183
$ws.Range("K51").Value = 179
$ws.Range("K52").Value = 389
$ws.Range("K53").Value = 193
$ws.Range("K54").Value = 265
$ws.Range("K55").Value = 162
$ws.Range("K57").Value = 50
$ws.Range("K60").Value = 92
$ws.Range("K63").Value = 47
$ws.Range("K64").Value = 87
$ws.Range("K65").Value = 330
$ws.Range("K67").Value = 554
$ws.Range("K73").Value = 129
$ws.Range("K77").Value = 101
$ws.Range("K78").Value = 170
$ws.Range("K83").Value = 304
$ws.Range("K85").Value = 643
$ws.Range("K88").Value = 165
$ws.Range("K92").Value = 52
$ws.Range("K94").Value = 180
$ws.Range("K96").Value = 160
$ws.Range("K98").Value = 74
$ws.Range("K101").Value = 14308

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K2").Value = 162
$ws.Range("K3").Value = 191
$ws.Range("K4").Value = 28
$ws.Range("K7").Value = 554

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("K6").Value = 131
$ws.Range("K7").Value = 265

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K2").Value = 213
$ws.Range("K4").Value = 39
$ws.Range("K5").Value = 23
$ws.Range("K6").Value = 208
$ws.Range("K7").Value = 751

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("K3").Value = 43
$ws.Range("K4").Value = 26
$ws.Range("K7").Value = 183

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("K3").Value = 38
$ws.Range("K7").Value = 131

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("K2").Value = 30
$ws.Range("K6").Value = 28
$ws.Range("K7").Value = 80

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K2").Value = 142
$ws.Range("K6").Value = 190
$ws.Range("K7").Value = 521

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("K4").Value = 16
$ws.Range("K7").Value = 170

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("K3").Value = 44
$ws.Range("K7").Value = 162

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range("K3").Value = 8
$ws.Range("K7").Value = 33

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("K6").Value = 35
$ws.Range("K7").Value = 146

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("K2").Value = 50
$ws.Range("K7").Value = 160

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("K2").Value = 16
$ws.Range("K7").Value = 87

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("K2").Value = 110
$ws.Range("K5").Value = 5
$ws.Range("K6").Value = 98
$ws.Range("K7").Value = 321

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("K2").Value = 71
$ws.Range("K6").Value = 41
$ws.Range("K7").Value = 181

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K2").Value = 150
$ws.Range("K3").Value = 137
$ws.Range("K4").Value = 15
$ws.Range("K7").Value = 420

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("K6").Value = 79
$ws.Range("K7").Value = 180

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("K2").Value = 23
$ws.Range("K7").Value = 66

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("K3").Value = 40
$ws.Range("K7").Value = 147

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("K3").Value = 14
$ws.Range("K7").Value = 74

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("K4").Value = 16
$ws.Range("K7").Value = 280

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("K6").Value = 50
$ws.Range("K7").Value = 129

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range("K6").Value = 25
$ws.Range("K7").Value = 52

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("K2").Value = 37
$ws.Range("K3").Value = 53
$ws.Range("K7").Value = 165

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("K2").Value = 49
$ws.Range("K6").Value = 59
$ws.Range("K7").Value = 179

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range("K3").Value = 10
$ws.Range("K7").Value = 50

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("K2").Value = 32
$ws.Range("K6").Value = 24
$ws.Range("K7").Value = 92

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K2").Value = 227
$ws.Range("K3").Value = 214
$ws.Range("K6").Value = 150
$ws.Range("K7").Value = 643

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("K3").Value = 36
$ws.Range("K7").Value = 101

$ws = $wb.Worksheets.Item('Jackson Park')
$ws.Range("K2").Value = 3
$ws.Range("K7").Value = 17

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("K3").Value = 101
$ws.Range("K4").Value = 23
$ws.Range("K7").Value = 389

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range("K3").Value = 3
$ws.Range("K7").Value = 43
